$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits --------------------------------------------------------
# "Ma hoc vien" / "ID hoc vien" keep the same digits but switch from a
# stored number to stored text (the cells already carry the Text number
# format, s="2").
$ws.Range("A2").Value = "1234567890"
$ws.Range("B2").Value = "1234567890"

# --- Row 3: a brand-new learner record -----------------------------------
# New cells default to the General number format, which would eat the
# leading zero on the phone number and turn everything into real numbers,
# so force Text formatting first.
$ws.Range("A3:L3").NumberFormat = "@"
$ws.Range("A3").Value = "1234567891"
$ws.Range("B3").Value = "1234567891"
$ws.Range("C3").Value = "1"
$ws.Range("E3").Value = "0123456790"
$ws.Range("F3").Value = "1234567890"
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "1"
$ws.Range("K3").Value = "Hà Nội"
$ws.Range("L3").Value = "Hà Nội"

# --- Column D (e-mail + hyperlink) --------------------------------------
# Stash D2's current Hyperlink look in a scratch cell - Hyperlinks.Add
# always reformats the target range, so we restore the real formatting
# once the links are in place.
$ws.Range("D2").Copy()
$ws.Range("Z1").PasteSpecial(-4122) # xlPasteFormats

$ws.Hyperlinks.Delete()
$ws.Range("D2").Value = "test2@gmail.com"
$ws.Range("D3").Value = "test3@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:test2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:test3@gmail.com")

$ws.Range("Z1").Copy()
$ws.Range("D2:D3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("Z1").Clear()
$excel.CutCopyMode = $false

# --- Selection ------------------------------------------------------------
$ws.Range("B4").Select()
